# Auto-generated: apply price/volume updates from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. A leading apostrophe forces the
# cell to stay a text value (matching the inline-string cells already in
# the sheet) instead of being auto-converted to a number/percentage.
$updates = @{
    'D2' = "'256.45"
    'E2' = "'0.42%"
    'D3' = "'26.97"
    'E3' = "'-4.44%"
    'D4' = "'4.763"
    'E4' = "'-9.22%"
    'D5' = "'0.05930"
    'E5' = "'1.20%"
    'E6' = "'-0.99%"
    'D7' = "'0.8670"
    'E7' = "'0.10%"
    'D8' = "'0.9398"
    'E8' = "'-6.27%"
    'D9' = "'0.1399"
    'E9' = "'-0.70%"
    'D10' = "'0.03740"
    'E10' = "'7.34%"
    'D11' = "'0.07121"
    'E11' = "'-0.55%"
    'D12' = "'0.03166"
    'E12' = "'-0.45%"
    'D13' = "'0.09240"
    'E13' = "'0.22%"
    'D14' = "'0.001550"
    'E14' = "'-0.04%"
    'D15' = "'0.006108"
    'E15' = "'5.06%"
    'D16' = "'3.482"
    'E16' = "'-0.51%"
    'D17' = "'3.197"
    'E17' = "'-0.87%"
    'E18' = "'0.69%"
    'D19' = "'0.01044"
    'E19' = "'1,627.56%"
    'D20' = "'0.3159"
    'E20' = "'-0.63%"
    'D21' = "'0.1298"
    'D22' = "'3.804"
    'E22' = "'7.64%"
    'D23' = "'0.04214"
    'E23' = "'1.40%"
    'E24' = "'-0.05%"
    'D25' = "'0.001224"
    'E25' = "'-0.14%"
    'D26' = "'0.004492"
    'E26' = "'-6.53%"
    'D27' = "'0.0001199"
    'E27' = "'-0.12%"
    'D28' = "'0.0001494"
    'E28' = "'1.86%"
    'D40' = "'0.03836"
    'E40' = "'0.80%"
    'D41' = "'0.006087"
    'E41' = "'3.99%"
    'D42' = "'0.1101"
    'E42' = "'0.07%"
    'D43' = "'0.002250"
    'E43' = "'-3.80%"
    'D44' = "'0.01110"
    'E44' = "'14.62%"
    'D45' = "'0.00005497"
    'E45' = "'5.16%"
    'E46' = "'-0.12%"
    'D47' = "'0.08847"
    'E47' = "'-4.92%"
    'D48' = "'0.002407"
    'E48' = "'11.86%"
    'D49' = "'0.00002099"
    'E49' = "'-0.12%"
    'D50' = "'0.0001999"
    'E50' = "'-0.12%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
